$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24, shifting existing rows 24-132 down to 25-133
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with data (a new "Espárragos" price record)
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = 'Macroferia Regional de Talca'
$ws.Range("C24").Value = 'Maule'
$ws.Range("D24").Value = 45243
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 300000000
$ws.Range("G24").Value = 'Espárragos'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 1400
$ws.Range("L24").Value = 1400
$ws.Range("M24").Value = 1400
$ws.Range("N24").Value = '$/kilo'
$ws.Range("O24").Value = 'Provincia de Linares'
$ws.Range("P24").Value = 1400
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 'Hortaliza'
